# Updated cryptos list values (price + 1h volume change) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.401.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.915.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  +0.74%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4796"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4037"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08185"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.005"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.34"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.921.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.029"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.208"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06834"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.012"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001037"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.009"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.410.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.651"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.188"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.114.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.562"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.093"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.013"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09596"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.592"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.550"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.365"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06327"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02280"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.179"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5917"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("E40").Value = "  +4.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.893"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1841"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.281"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.395"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07465"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5536"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.920"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.422"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.24%  "
